$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds the "Förändrad" (changed) date as a serial number.
# Every data row (2 through 97) had its value changed from 45190 to 45192.
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45190) {
        $cell.Value2 = 45192
    }
}
